$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.389.06"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "3.488.51"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'610.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.09%  "
$ws.Range("D6").Value = "'186.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.214"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'53.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "'9.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "4.033.96"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "'601.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").Value = "69.462.13"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "3.492.33"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'0.985"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").Value = "'17.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").Value = "'105.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.34%  "
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "'5.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'3.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("D29").Value = "'33.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "'4.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +21.99%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'12.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "'3.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'523.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.45%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.613.83"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.396"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "'36.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.0458"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("D48").Value = "'8.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("D51").Value = "'1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.09%  "
